$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tile_spat_count_entry")

for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 2).Value = 43468
    $ws.Cells.Item($r, 3).Value = 43497
}

$ws.Range("A2:P28").Select()
